$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row containing "syntok" in column A and delete the entire row
$found = $ws.Range("A1:A38").Find("syntok")
if ($found -ne $null) {
    $rowNum = $found.Row
    $ws.Rows.Item($rowNum).Delete()
}
